$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new day's entry (row 11) that was previously blank except for the date.
$ws.Range("D11").Value = 4385
$ws.Range("F11").Value = 0
$ws.Range("H11").Value = 329
$ws.Range("J11").Value = 50
$ws.Range("L11").Value = "Return Trip by Electric Wizard"

# Move the active selection to reflect where editing continued (H12).
[void]$ws.Range("H12").Select()
